# Implement first version of lot sizing rules
# - Bump NrBuckets (Generic!B4) from 3 to 4
# - Bump Productdata minimum lot size column (C4:C6) from 5 to 10
# - Extend ForecastedAverageDemand and ForcastedStandardDeviation with a new
#   4th bucket (row 5) matching the pattern of the existing rows.

$wb = $excel.ActiveWorkbook

# --- Generic sheet: NrBuckets 3 -> 4 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 4

# --- Productdata sheet: lot size 5 -> 10 for rows 4-6 ---
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("C4").Value = 10
$wsProduct.Range("C5").Value = 10
$wsProduct.Range("C6").Value = 10

# --- ForecastedAverageDemand sheet: add new bucket row 5 ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("A4:I4").Copy()
$wsAvgDemand.Range("A5:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsAvgDemand.Range("A5").Value = 3
$wsAvgDemand.Range("B5").Value = 0
$wsAvgDemand.Range("C5").Value = 0
$wsAvgDemand.Range("D5").Value = 0
$wsAvgDemand.Range("E5").Value = 0
$wsAvgDemand.Range("F5").Value = 0
$wsAvgDemand.Range("G5").Value = 1
$wsAvgDemand.Range("H5").Value = 1
$wsAvgDemand.Range("I5").Value = 1

# --- ForcastedStandardDeviation sheet: add new bucket row 5 ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("A4:I4").Copy()
$wsStdDev.Range("A5:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsStdDev.Range("A5").Value = 3
$wsStdDev.Range("B5").Value = 0
$wsStdDev.Range("C5").Value = 0
$wsStdDev.Range("D5").Value = 0
$wsStdDev.Range("E5").Value = 0
$wsStdDev.Range("F5").Value = 0
$wsStdDev.Range("G5").Value = 1
$wsStdDev.Range("H5").Value = 1
$wsStdDev.Range("I5").Value = 1
